# Auto-generated edit script applying the crypto price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a Number by Excel's
# type coercion (plain decimal strings like "207.60") are first forced to the
# Text number format so the assigned value stays a literal string (matching the
# original inlineStr cells), then the format is reset back to the workbook's
# default style so no stray formatting is introduced.
$numericLooking = @('D5', 'D6', 'D8', 'D10', 'D11', 'D14', 'D15', 'D16', 'D19', 'D20', 'D22', 'D23', 'D24', 'D25', 'D26', 'D27', 'D30', 'D32', 'D33', 'D39', 'D40', 'D41', 'D45', 'D48', 'D49', 'D51')
foreach ($addr in $numericLooking) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '26.953.56'
$ws.Range('E2').Value = '  -0.35%  '
$ws.Range('D3').Value = '1.563.10'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '207.60'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').Value = '0.490'
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').Value = '22.10'
$ws.Range('E8').Value = '  +0.42%  '
$ws.Range('E9').Value = '  -0.70%  '
$ws.Range('D10').Value = '0.0600'
$ws.Range('E10').Value = '  +1.99%  '
$ws.Range('D11').Value = '0.0857'
$ws.Range('E11').Value = '  -0.27%  '
$ws.Range('D12').Value = '1.784.75'
$ws.Range('E12').Value = '  -0.31%  '
$ws.Range('D13').Value = '1.564.89'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('D14').Value = '3.76'
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').Value = '0.520'
$ws.Range('E15').Value = '  -0.20%  '
$ws.Range('D16').Value = '62.11'
$ws.Range('E16').Value = '  +0.17%  '
$ws.Range('D17').Value = '26.946.97'
$ws.Range('E17').Value = '  -0.35%  '
$ws.Range('D18').Value = '0.0₃0706'
$ws.Range('E18').Value = '  +1.20%  '
$ws.Range('D19').Value = '217.05'
$ws.Range('E19').Value = '  -0.62%  '
$ws.Range('D20').Value = '7.36'
$ws.Range('E20').Value = '  +0.49%  '
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').Value = '4.11'
$ws.Range('E22').Value = '  +0.84%  '
$ws.Range('D23').Value = '9.23'
$ws.Range('E23').Value = '  -0.65%  '
$ws.Range('D24').Value = '1.94'
$ws.Range('E24').Value = '  -1.06%  '
$ws.Range('D25').Value = '152.37'
$ws.Range('E25').Value = '  -1.28%  '
$ws.Range('D26').Value = '6.61'
$ws.Range('E26').Value = '  +0.00%  '
$ws.Range('D27').Value = '15.09'
$ws.Range('E27').Value = '  +0.80%  '
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('E29').Value = '  -0.17%  '
$ws.Range('D30').Value = '0.0472'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  +1.16%  '
$ws.Range('D32').Value = '3.24'
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').Value = '3.12'
$ws.Range('E33').Value = '  +2.16%  '
$ws.Range('D34').Value = '1.421.26'
$ws.Range('E34').Value = '  -1.78%  '
$ws.Range('E35').Value = '  +2.96%  '
$ws.Range('E36').Value = '  +11.54%  '
$ws.Range('E37').Value = '  +1.59%  '
$ws.Range('E38').Value = '  +0.07%  '
$ws.Range('D39').Value = '0.533'
$ws.Range('E39').Value = '  +1.88%  '
$ws.Range('D40').Value = '0.809'
$ws.Range('E40').Value = '  -0.82%  '
$ws.Range('D41').Value = '5.78'
$ws.Range('E41').Value = '  +0.11%  '
$ws.Range('E42').Value = '  -0.15%  '
$ws.Range('E43').Value = '  +2.03%  '
$ws.Range('E44').Value = '  +1.50%  '
$ws.Range('D45').Value = '64.89'
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('E46').Value = '  -0.76%  '
$ws.Range('D47').Value = '1.697.95'
$ws.Range('E47').Value = '  -0.36%  '
$ws.Range('D48').Value = '87.62'
$ws.Range('E48').Value = '  +1.17%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '0.0521'
$ws.Range('E49').Value = '  -0.27%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0102'
$ws.Range('E50').Value = '  +8.80%  '
$ws.Range('D51').Value = '0.0961'
$ws.Range('E51').Value = '  -0.58%  '

foreach ($addr in $numericLooking) {
    $ws.Range($addr).Style = "Normal"
}
